# Refresh the demo user data from the "23" batch to the "24" batch.
#
# - "login" sheet (G2:I21): hard-coded user name / user name / email values
#   carry a "23" suffix; bump it to "24".
# - "Sheet1" sheet: I2 is the literal suffix number (23) that the
#   CONCATENATE formulas in A2:C21 (and A3:B21 / C3:C21 shared formulas)
#   pull from; bumping it to 24 ripples the new suffix through those rows.
#   Its selection also moves back up to the refreshed A2:C21 block.
# - "order" sheet data (R2:T21, the "28" batch) is untouched.
# - Active tab moves from "order" back to "login".

$wb = $excel.ActiveWorkbook

$login = $wb.Worksheets.Item("login")
$order = $wb.Worksheets.Item("order")
$sheet1 = $wb.Worksheets.Item("Sheet1")

# --- "login": bump the displayed "23" user batch to "24" ---
for ($r = 2; $r -le 21; $r++) {
    $g = $login.Cells.Item($r, 7).Value()
    $h = $login.Cells.Item($r, 8).Value()
    $i = $login.Cells.Item($r, 9).Value()

    $login.Cells.Item($r, 7).Value = ($g -replace '23', '24')
    $login.Cells.Item($r, 8).Value = ($h -replace '23', '24')
    $login.Cells.Item($r, 9).Value = ($i -replace '23', '24')
}

# --- "Sheet1": bump the suffix driving the CONCATENATE formulas ---
$sheet1.Cells.Item(2, 9).Value = 24

# Refresh Sheet1's own selection to the just-updated block and drop the
# old scroll/selection that pointed at the "28" block further down.
$sheet1.Activate()
$sheet1.Range("A2:C21").Select()

# "order" keeps its own selection as-is (still R2:T21 / "28" data).
$order.Activate()
$order.Range("R2:T21").Select()

# "login" becomes the active/selected tab again.
$login.Activate()
$login.Range("G2:I21").Select()
